$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 159 (pushes existing rows 159:172 down to 160:173,
# copying the formatting of row 159 - notably the date style on column D).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(159, 1).Value  = 5
$ws.Cells.Item(159, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(159, 3).Value  = "Maule"
$ws.Cells.Item(159, 4).Value  = 44461
$ws.Cells.Item(159, 5).Value  = 7
$ws.Cells.Item(159, 6).Value  = 100114014
$ws.Cells.Item(159, 7).Value  = "Betarraga"
$ws.Cells.Item(159, 8).Value  = "Sin especificar"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 3000
$ws.Cells.Item(159, 11).Value = 700
$ws.Cells.Item(159, 12).Value = 700
$ws.Cells.Item(159, 13).Value = 700
$ws.Cells.Item(159, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(159, 15).Value = "Región del Maule"
$ws.Cells.Item(159, 16).Value = 140
$ws.Cells.Item(159, 17).Value = 5
$ws.Cells.Item(159, 18).Value = "Hortaliza"
